# Node and Html dashboard code
# Add two more link rows (esign.verasys.in, vsign.in) below the existing
# list of links, turning them into real hyperlinks, then resize column A
# to fit the new (longer) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: https://esign.verasys.in/
$ws.Hyperlinks.Add($ws.Range("A7"), "https://esign.verasys.in/")

# Row 8: https://vsign.in/
$ws.Hyperlinks.Add($ws.Range("A8"), "https://vsign.in/")

# Keep the selection on the last entered cell, like Excel leaves it
# after typing/tabbing through the new rows.
[void]$ws.Range("A8").Select()

# Column A now holds longer URLs -- best-fit the width like Excel would
# after double-clicking the column border.
$ws.Columns("A:A").AutoFit()
